{"js": "// Apply the LOQ4252.docx content rotation described by the diff:\n//  - \"Objetivos\" body        <- old \"Programa resumido\" body text\n//  - \"Docente(s) ...\" body   <- old \"Objetivos\" body text\n//  - \"Programa resumido\" body<- old \"Programa\" body text (big numbered list)\n//  - \"Programa\" body         <- old \"Avalia\u00e7\u00e3o/M\u00e9todo\" content\n//  - \"Avalia\u00e7\u00e3o\" (M\u00e9todo/Crit\u00e9rio/Norma) content runs each shift to the next label,\n//    and \"Norma de recupera\u00e7\u00e3o\" gets the old \"Bibliografia\" body text\n//  - \"Bibliografia\" body     <- old \"Docente(s) ...\" body text (\"4808662 - ...\")\n//\n// All paragraph styles/structure stay the same; only the w:t / w:br content inside\n// the (non-heading, non-bold-label) runs is replaced.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text,style\");\nawait context.sync();\n\n// Helper: find the paragraph that immediately follows the (unique) heading paragraph\n// whose text equals `headingText`.\nfunction bodyParagraphAfterHeading(headingText) {\n  const items = paragraphs.items;\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === headingText) {\n      return items[i + 1];\n    }\n  }\n  throw new Error(\"Heading not found: \" + headingText);\n}\n\nconst objetivosBody = bodyParagraphAfterHeading(\"Objetivos\");\nconst docentesBody = bodyParagraphAfterHeading(\"Docente(s) Respons\u00e1vel(eis) \");\nconst programaResumidoBody = bodyParagraphAfterHeading(\"Programa resumido\");\nconst programaBody = bodyParagraphAfterHeading(\"Programa\");\nconst avaliacaoBody = bodyParagraphAfterHeading(\"Avalia\u00e7\u00e3o\");\nconst bibliografiaBody = bodyParagraphAfterHeading(\"Bibliografia\");\n\n// New text values (\\v == vertical tab == Word's manual line break == <w:br/>).\nconst textObjetivos =\n  \"Conceitos ligados ao escoamento de flu\u00eddos e equa\u00e7\u00f5es fundamentais, Escoamento incompress\u00edvel de fluidos n\u00e3o viscosos, Escoamento viscoso incompress\u00edvel, Transfer\u00eancia de Calor. Transfer\u00eancia de Massa\";\n\nconst textDocentes =\n  \"Fornecer os conceitos b\u00e1sicos de Mec\u00e2nica dos Fluidos e Transfer\u00eancia de Calor e Massa com aplica\u00e7\u00f5es \u00e0 Engenharia. Possibilitar aos alunos uma base cient\u00edfica para que possam se desenvolver em demais disciplinas tecnol\u00f3gicas do curso.\";\n\nconst textProgramaResumido =\n  \"1. Conceitos ligados ao escoamento de flu\u00eddos e equa\u00e7\u00f5es fundamentais\\v\" +\n  \"1.1. Caracter\u00edsticas e defini\u00e7\u00f5es dos escoamentos;\\v\" +\n  \"1.2. Conceitos de sistema e volume de controle;\\v\" +\n  \"1.3. Equa\u00e7\u00e3o da conserva\u00e7\u00e3o da massa;\\v\" +\n  \"1.4. Equa\u00e7\u00e3o da conserva\u00e7\u00e3o da energia;\\v\" +\n  \"1.5. Equa\u00e7\u00e3o da conserva\u00e7\u00e3o da quantidade de movimento;\\v\" +\n  \"1.6. Introdu\u00e7\u00e3o \u00e0 an\u00e1lise diferencial do movimento de fluidos.\\v\" +\n  \"\\v\" +\n  \"2. Escoamento incompress\u00edvel de fluidos n\u00e3o viscosos\\v\" +\n  \"2.1. Equa\u00e7\u00e3o de Euler;\\v\" +\n  \"2.2. Equa\u00e7\u00e3o de Bernoulli;\\v\" +\n  \"2.3. Aplica\u00e7\u00f5es da equa\u00e7\u00e3o de Bernoulli.\\v\" +\n  \"\\v\" +\n  \"3. Escoamento viscoso incompress\u00edvel\\v\" +\n  \"3.1. Atrito e perda de carga;\\v\" +\n  \"3.2. Avalia\u00e7\u00e3o das perdas de carga: regime laminar e turbulento;\\v\" +\n  \"3.3. Equa\u00e7\u00f5es de Hagen - Poiseuille e Darcy \u2013 Weisbach\\v\" +\n  \"3.4. Diagrama de Moody e Moody \u2013Rouse;\\v\" +\n  \"3.5. M\u00e9todo dos comprimentos equivalentes.\\v\" +\n  \"3.6. Presen\u00e7a de m\u00e1quina no escoamento (bomba e turbina), Pot\u00eancia e rendimento;\\v\" +\n  \"3.7. Medidores de vaz\u00e3o.\\v\" +\n  \"\\v\" +\n  \"4. Transfer\u00eancia de Calor\\v\" +\n  \"4.1. Defini\u00e7\u00e3o de Calor.\\v\" +\n  \"4.2. Mecanismo da Condu\u00e7\u00e3o.\\v\" +\n  \"4.3. Mecanismo da Convec\u00e7\u00e3o.\\v\" +\n  \"4.4. Associa\u00e7\u00e3o de Mecanismos.\\v\" +\n  \"\\v\" +\n  \"5. Transfer\u00eancia de Massa\\v\" +\n  \"5.1. Difus\u00e3o e convec\u00e7\u00e3o m\u00e1ssica;\\v\" +\n  \"5.2. 1\u00aa lei de Fick;\\v\" +\n  \"5.3. Concentra\u00e7\u00f5es m\u00e1ssica e molar;\\v\" +\n  \"5.4. Fra\u00e7\u00f5es m\u00e1ssica e molar;\\v\" +\n  \"5.5. Velocidades m\u00e9dias m\u00e1ssica e molar;\\v\" +\n  \"5.6. Fluxos difusivo m\u00e1ssico, difusivo molar, convectivo m\u00e1ssico e convectivo molar;\\v\" +\n  \"5.7. Fluxo m\u00e1ssico total e fluxo molar total.\";\n\nconst textPrograma =\n  \"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de exerc\u00edcios, aulas de laborat\u00f3rio.\";\n\nconst textMetodoContent =\n  \"Nota de duas provas (P1 e P2)\\vF\u00f3rmula: M1 = (P1 + 2 x P2)/3..\\v\";\n\nconst textCriterioContent =\n  \"Aplica\u00e7\u00e3o de uma prova envolvendo o assunto de todo semestre.\\vNR (nota da recupera\u00e7\u00e3o) = (M1 + NR)/2.\\v\";\n\nconst textNormaContent =\n  \"1. FOX, R.W., MCDONALD, A.T., \\u201cIntrodu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos Fluidos\\u201d, Ed. Guanabara Koogan.\\v\" +\n  \"2. STREETER, V.L., WYLE,E.B., \\u201cMec\u00e2nica dos Fluidos\\u201d, Ed. Mc Graw Hill.\\v\" +\n  \"3. OZISIK,M.N., \\u201cTransfer\u00eancia de Calor.\\u201d, Ed. Guanabara Koogan.\\v\" +\n  \"4. INCROPERA, F.P.W., \\u201cFundamentos de Transfer\u00eancia de Calor e Massa\\u201d, Ed. Guanabara Koogan.\\v\" +\n  \"5. MUNSON, B.R.; YOUNG, D.F.; OKIISHI, T.H. Fundamentos da Mec\u00e2nica dos Fluidos. Editora Edgard Blucher\\v\" +\n  \"6 - GIORGETI, M. (2012) Fundamentos de Fen\u00f4menos de Transporte. Editora Campus\";\n\nconst textBibliografia = \"4808662 - Lucr\u00e9cio F\u00e1bio dos Santos\";\n\n// Simple single-run body paragraphs: replace the whole paragraph text in place\n// (this preserves the paragraph mark / style, only swaps its run content).\nobjetivosBody.insertText(textObjetivos, Word.InsertLocation.replace);\ndocentesBody.insertText(textDocentes, Word.InsertLocation.replace);\nprogramaResumidoBody.insertText(textProgramaResumido, Word.InsertLocation.replace);\nprogramaBody.insertText(textPrograma, Word.InsertLocation.replace);\nbibliografiaBody.insertText(textBibliografia, Word.InsertLocation.replace);\nawait context.sync();\n\n// \"Avalia\u00e7\u00e3o\" paragraph has 3 bold labels (\"M\u00e9todo: \", \"Crit\u00e9rio: \",\n// \"Norma de recupera\u00e7\u00e3o: \") each followed by a non-bold content run. Isolate each\n// content span using the (stable) labels as anchors, so the bold runs are untouched.\nconst avalRange = avaliacaoBody.getRange();\n\nconst metodoLabel = avalRange.search(\"M\u00e9todo: \", { matchCase: true });\nconst criterioLabel = avalRange.search(\"Crit\u00e9rio: \", { matchCase: true });\nconst normaLabel = avalRange.search(\"Norma de recupera\u00e7\u00e3o: \", { matchCase: true });\nmetodoLabel.load(\"items\");\ncriterioLabel.load(\"items\");\nnormaLabel.load(\"items\");\nawait context.sync();\n\nconst metodoContentRange = metodoLabel.items[0]\n  .getRange(\"After\")\n  .expandTo(criterioLabel.items[0].getRange(\"Before\"));\nconst criterioContentRange = criterioLabel.items[0]\n  .getRange(\"After\")\n  .expandTo(normaLabel.items[0].getRange(\"Before\"));\nconst normaContentRange = normaLabel.items[0]\n  .getRange(\"After\")\n  .expandTo(avaliacaoBody.getRange(\"End\"));\n\nmetodoContentRange.insertText(textMetodoContent, Word.InsertLocation.replace);\ncriterioContentRange.insertText(textCriterioContent, Word.InsertLocation.replace);\nnormaContentRange.insertText(textNormaContent, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Apply the LOQ4252.docx content rotation described by the diff:\n#  - \"Objetivos\" body        <- old \"Programa resumido\" body text\n#  - \"Docente(s) ...\" body   <- old \"Objetivos\" body text\n#  - \"Programa resumido\" body<- old \"Programa\" body text (big numbered list)\n#  - \"Programa\" body         <- old \"Avalia\u00e7\u00e3o/M\u00e9todo\" content\n#  - \"Avalia\u00e7\u00e3o\" (M\u00e9todo/Crit\u00e9rio/Norma) content runs each shift to the next label,\n#    and \"Norma de recupera\u00e7\u00e3o\" gets the old \"Bibliografia\" body text\n#  - \"Bibliografia\" body     <- old \"Docente(s) ...\" body text (\"4808662 - ...\")\n#\n# All paragraph styles/structure stay the same; only the w:t / w:br content inside\n# the (non-heading, non-bold-label) runs is replaced.\n\n$d = $word.ActiveDocument\n\nfunction Get-BodyParagraphAfterHeading {\n    param($doc, [string]$HeadingText)\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $txt = $p.Range.Text.TrimEnd([char]13)\n        if ($txt -eq $HeadingText) {\n            return $doc.Paragraphs.Item($i + 1)\n        }\n    }\n    throw \"Heading not found: $HeadingText\"\n}\n\n$objetivosBody = Get-BodyParagraphAfterHeading $d \"Objetivos\"\n$docentesBody = Get-BodyParagraphAfterHeading $d \"Docente(s) Respons\u00e1vel(eis) \"\n$programaResumidoBody = Get-BodyParagraphAfterHeading $d \"Programa resumido\"\n$programaBody = Get-BodyParagraphAfterHeading $d \"Programa\"\n$avaliacaoBody = Get-BodyParagraphAfterHeading $d \"Avalia\u00e7\u00e3o\"\n$bibliografiaBody = Get-BodyParagraphAfterHeading $d \"Bibliografia\"\n\n$NL = [char]11   # Word's manual line break char -> becomes <w:br/>\n\n$textObjetivos = \"Conceitos ligados ao escoamento de flu\u00eddos e equa\u00e7\u00f5es fundamentais, Escoamento incompress\u00edvel de fluidos n\u00e3o viscosos, Escoamento viscoso incompress\u00edvel, Transfer\u00eancia de Calor. Transfer\u00eancia de Massa\"\n\n$textDocentes = \"Fornecer os conceitos b\u00e1sicos de Mec\u00e2nica dos Fluidos e Transfer\u00eancia de Calor e Massa com aplica\u00e7\u00f5es \u00e0 Engenharia. Possibilitar aos alunos uma base cient\u00edfica para que possam se desenvolver em demais disciplinas tecnol\u00f3gicas do curso.\"\n\n$textProgramaResumido = (\n    \"1. Conceitos ligados ao escoamento de flu\u00eddos e equa\u00e7\u00f5es fundamentais\" + $NL +\n    \"1.1. Caracter\u00edsticas e defini\u00e7\u00f5es dos escoamentos;\" + $NL +\n    \"1.2. Conceitos de sistema e volume de controle;\" + $NL +\n    \"1.3. Equa\u00e7\u00e3o da conserva\u00e7\u00e3o da massa;\" + $NL +\n    \"1.4. Equa\u00e7\u00e3o da conserva\u00e7\u00e3o da energia;\" + $NL +\n    \"1.5. Equa\u00e7\u00e3o da conserva\u00e7\u00e3o da quantidade de movimento;\" + $NL +\n    \"1.6. Introdu\u00e7\u00e3o \u00e0 an\u00e1lise diferencial do movimento de fluidos.\" + $NL +\n    $NL +\n    \"2. Escoamento incompress\u00edvel de fluidos n\u00e3o viscosos\" + $NL +\n    \"2.1. Equa\u00e7\u00e3o de Euler;\" + $NL +\n    \"2.2. Equa\u00e7\u00e3o de Bernoulli;\" + $NL +\n    \"2.3. Aplica\u00e7\u00f5es da equa\u00e7\u00e3o de Bernoulli.\" + $NL +\n    $NL +\n    \"3. Escoamento viscoso incompress\u00edvel\" + $NL +\n    \"3.1. Atrito e perda de carga;\" + $NL +\n    \"3.2. Avalia\u00e7\u00e3o das perdas de carga: regime laminar e turbulento;\" + $NL +\n    \"3.3. Equa\u00e7\u00f5es de Hagen - Poiseuille e Darcy \u2013 Weisbach\" + $NL +\n    \"3.4. Diagrama de Moody e Moody \u2013Rouse;\" + $NL +\n    \"3.5. M\u00e9todo dos comprimentos equivalentes.\" + $NL +\n    \"3.6. Presen\u00e7a de m\u00e1quina no escoamento (bomba e turbina), Pot\u00eancia e rendimento;\" + $NL +\n    \"3.7. Medidores de vaz\u00e3o.\" + $NL +\n    $NL +\n    \"4. Transfer\u00eancia de Calor\" + $NL +\n    \"4.1. Defini\u00e7\u00e3o de Calor.\" + $NL +\n    \"4.2. Mecanismo da Condu\u00e7\u00e3o.\" + $NL +\n    \"4.3. Mecanismo da Convec\u00e7\u00e3o.\" + $NL +\n    \"4.4. Associa\u00e7\u00e3o de Mecanismos.\" + $NL +\n    $NL +\n    \"5. Transfer\u00eancia de Massa\" + $NL +\n    \"5.1. Difus\u00e3o e convec\u00e7\u00e3o m\u00e1ssica;\" + $NL +\n    \"5.2. 1\u00aa lei de Fick;\" + $NL +\n    \"5.3. Concentra\u00e7\u00f5es m\u00e1ssica e molar;\" + $NL +\n    \"5.4. Fra\u00e7\u00f5es m\u00e1ssica e molar;\" + $NL +\n    \"5.5. Velocidades m\u00e9dias m\u00e1ssica e molar;\" + $NL +\n    \"5.6. Fluxos difusivo m\u00e1ssico, difusivo molar, convectivo m\u00e1ssico e convectivo molar;\" + $NL +\n    \"5.7. Fluxo m\u00e1ssico total e fluxo molar total.\"\n)\n\n$textPrograma = \"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de exerc\u00edcios, aulas de laborat\u00f3rio.\"\n\n$textMetodoContent = \"Nota de duas provas (P1 e P2)\" + $NL + \"F\u00f3rmula: M1 = (P1 + 2 x P2)/3..\" + $NL\n\n$textCriterioContent = \"Aplica\u00e7\u00e3o de uma prova envolvendo o assunto de todo semestre.\" + $NL + \"NR (nota da recupera\u00e7\u00e3o) = (M1 + NR)/2.\" + $NL\n\n$textNormaContent = (\n    \"1. FOX, R.W., MCDONALD, A.T., \u201cIntrodu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos Fluidos\u201d, Ed. Guanabara Koogan.\" + $NL +\n    \"2. STREETER, V.L., WYLE,E.B., \u201cMec\u00e2nica dos Fluidos\u201d, Ed. Mc Graw Hill.\" + $NL +\n    \"3. OZISIK,M.N., \u201cTransfer\u00eancia de Calor.\u201d, Ed. Guanabara Koogan.\" + $NL +\n    \"4. INCROPERA, F.P.W., \u201cFundamentos de Transfer\u00eancia de Calor e Massa\u201d, Ed. Guanabara Koogan.\" + $NL +\n    \"5. MUNSON, B.R.; YOUNG, D.F.; OKIISHI, T.H. Fundamentos da Mec\u00e2nica dos Fluidos. Editora Edgard Blucher\" + $NL +\n    \"6 - GIORGETI, M. (2012) Fundamentos de Fen\u00f4menos de Transporte. Editora Campus\"\n)\n\n$textBibliografia = \"4808662 - Lucr\u00e9cio F\u00e1bio dos Santos\"\n\n# Simple single-run body paragraphs: replace the whole paragraph text in place\n# (Range.Text on a Paragraphs.Item().Range excludes the trailing pilcrow, so this\n# preserves the paragraph mark / style and only swaps its run content).\n$objetivosBody.Range.Text = $textObjetivos\n$docentesBody.Range.Text = $textDocentes\n$programaResumidoBody.Range.Text = $textProgramaResumido\n$programaBody.Range.Text = $textPrograma\n$bibliografiaBody.Range.Text = $textBibliografia\n\n# \"Avalia\u00e7\u00e3o\" paragraph has 3 bold labels (\"M\u00e9todo: \", \"Crit\u00e9rio: \",\n# \"Norma de recupera\u00e7\u00e3o: \") each followed by a non-bold content run. Isolate each\n# content span using the (stable) labels as anchors, so the bold runs are untouched.\n$avalRange = $avaliacaoBody.Range\n\n$findMetodo = $avalRange.Duplicate\n$findMetodo.Find.ClearFormatting()\n$findMetodo.Find.Text = \"M\u00e9todo: \"\n$findMetodo.Find.Execute() | Out-Null\n$afterMetodo = $findMetodo.End\n\n$findCriterio = $avalRange.Duplicate\n$findCriterio.Find.ClearFormatting()\n$findCriterio.Find.Text = \"Crit\u00e9rio: \"\n$findCriterio.Find.Execute() | Out-Null\n$beforeCriterio = $findCriterio.Start\n$afterCriterio = $findCriterio.End\n\n$findNorma = $avalRange.Duplicate\n$findNorma.Find.ClearFormatting()\n$findNorma.Find.Text = \"Norma de recupera\u00e7\u00e3o: \"\n$findNorma.Find.Execute() | Out-Null\n$beforeNorma = $findNorma.Start\n$afterNorma = $findNorma.End\n\n$endOfAvaliacao = $avaliacaoBody.Range.End\n\n# Replace right-to-left (Norma, then Crit\u00e9rio, then M\u00e9todo) so that the already-\n# computed offsets for the earlier (leftmost) ranges stay valid while later ones\n# are mutated first.\n$normaContentRange = $d.Range($afterNorma, $endOfAvaliacao)\n$normaContentRange.Text = $textNormaContent\n\n$criterioContentRange = $d.Range($afterCriterio, $beforeNorma)\n$criterioContentRange.Text = $textCriterioContent\n\n$metodoContentRange = $d.Range($afterMetodo, $beforeCriterio)\n$metodoContentRange.Text = $textMetodoContent\n\nWrite-Output \"done\"\n"}
